# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 2024-04-09 (45391) to 2024-04-11 (45393).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45391) {
        $cell.Value2 = 45393
    }
}
